# Weighted Evaluation Table for Project Selection - fill in the "Emre Deniz SENEL"
# reviewer block (columns I:O, rows 3-10) and fix the L13 score.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Software Based ---
$ws.Range("I3").Value = 5
$ws.Range("L3").Value = 9
$ws.Range("M3").Value = 7
$ws.Range("N3").Value = 7
$ws.Range("O3").Value = 9

# --- Row 4: Hardware Based ---
$ws.Range("I4").Value = 5
$ws.Range("L4").Value = 4
$ws.Range("M4").Value = 8
$ws.Range("N4").Value = 8
$ws.Range("O4").Value = 8

# --- Row 5: Multidisciplinary ---
$ws.Range("I5").Value = 5
$ws.Range("L5").Value = 2
$ws.Range("M5").Value = 9
$ws.Range("N5").Value = 9
$ws.Range("O5").Value = 6

# --- Row 6: Adoptability to Everyday Life ---
$ws.Range("I6").Value = 8
$ws.Range("L6").Value = 10
$ws.Range("M6").Value = 9
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 1

# --- Row 7: Mechanic Based ---
$ws.Range("I7").Value = 5
$ws.Range("L7").Value = 2
$ws.Range("M7").Value = 9
$ws.Range("N7").Value = 4
$ws.Range("O7").Value = 8

# --- Row 8: Open to Improvement ---
$ws.Range("I8").Value = 9
$ws.Range("L8").Value = 10
$ws.Range("M8").Value = 9
$ws.Range("N8").Value = 2
$ws.Range("O8").Value = 5

# --- Row 9: Enjoyable ---
$ws.Range("I9").Value = 4
$ws.Range("L9").Value = 6
$ws.Range("M9").Value = 8
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 6

# --- Row 10: Marketable ---
$ws.Range("I10").Value = 9
$ws.Range("L10").Value = 9
$ws.Range("M10").Value = 9
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = 7

# --- Row 13: score correction ---
$ws.Range("L13").Value = 0

# --- Selection state: the active cell/selection moves to N14 ---
$ws.Range("N14").Select()
